# Generate Report for Handback
#
# The source file "753516b8-1b0c-4f2a-a4fa-c621cde4422f.md" has been handed
# back (for both the zh-cn and de-de locales). Update the localization
# status report accordingly:
#   - Status moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" on every sheet that tracks it.
#   - The per-locale "Latest Handback DateTime" is stamped with the
#     handback time for each locale.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $handedBack
$overview.Range("C3").Value = $handedBack

# --- zh-cn sheet -----------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $handedBack
$zhcn.Range("H3").Value = "2016-03-31 05:24:55"

# --- de-de sheet -----------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $handedBack
$dede.Range("H3").Value = "2016-03-31 05:25:11"
